$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.806.38'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.042.17'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'227.49"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("D7").Value = "'59.67"
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '2.343.59'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = "'14.45"
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").Value = "'21.04"
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = "'5.48"
$ws.Range("E15").Value = '  +4.59%  '
$ws.Range("D16").Value = "'0.770"
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '2.047.62'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '37.756.26'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = "'69.51"
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = "'223.82"
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").Value = "'2.27"
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("D26").Value = "'169.46"
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").Value = "'18.79"
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("E32").Value = '  +9.17%  '
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("D34").Value = "'0.0602"
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = "'6.53"
$ws.Range("E36").Value = '  +1.57%  '
$ws.Range("E37").Value = '  +3.84%  '
$ws.Range("E38").Value = '  +6.15%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = "'18.01"
$ws.Range("E40").Value = '  +5.90%  '
$ws.Range("D41").Value = '1.527.59'
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = "'97.58"
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").Value = "'0.0908"
$ws.Range("D46").Value = "'4.16"
$ws.Range("E46").Value = '  +6.21%  '
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").Value = '2.232.94'
$ws.Range("E51").Value = '  +0.29%  '
